$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44511
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 900
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 950
$ws.Range("P2").Value = 950
$ws.Range("D4").Value = 44523
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 800
$ws.Range("L4").Value = 900
$ws.Range("M4").Value = 850
$ws.Range("P4").Value = 850
$ws.Range("D5").Value = 44524
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 800
$ws.Range("L5").Value = 900
$ws.Range("M5").Value = 850
$ws.Range("P5").Value = 850
$ws.Range("D6").Value = 44460
$ws.Range("H6").Value = "Verde"
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 2200
$ws.Range("L6").Value = 2300
$ws.Range("M6").Value = 2250
$ws.Range("P6").Value = 2250
$ws.Range("D7").Value = 44476
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 1100
$ws.Range("L7").Value = 1200
$ws.Range("M7").Value = 1150
$ws.Range("P7").Value = 1150
$ws.Range("D8").Value = 44512
$ws.Range("J8").Value = 600
$ws.Range("D9").Value = 44532
$ws.Range("J9").Value = 240
$ws.Range("D10").Value = 44508
$ws.Range("J10").Value = 400
$ws.Range("K10").Value = 900
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = 950
$ws.Range("P10").Value = 950
$ws.Range("D11").Value = 44530
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 800
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = 850
$ws.Range("P11").Value = 850
$ws.Range("D12").Value = 44545
$ws.Range("J12").Value = 4000
$ws.Range("D13").Value = 44837
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 1800
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = 1900
$ws.Range("P13").Value = 1900
$ws.Range("D14").Value = 44525
$ws.Range("J14").Value = 360
$ws.Range("D15").Value = 44827
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 2200
$ws.Range("L15").Value = 2300
$ws.Range("M15").Value = 2250
$ws.Range("P15").Value = 2250
$ws.Range("D16").Value = 44553
$ws.Range("J16").Value = 8000
$ws.Range("D17").Value = 44831
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 2000
$ws.Range("L17").Value = 2200
$ws.Range("M17").Value = 2100
$ws.Range("P17").Value = 2100
$ws.Range("D18").Value = 44516
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 900
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = 950
$ws.Range("P18").Value = 950
$ws.Range("D19").Value = 44510
$ws.Range("D20").Value = 44505
$ws.Range("J20").Value = 440
$ws.Range("D21").Value = 44848
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("J21").Value = 500
$ws.Range("K21").Value = 1300
$ws.Range("L21").Value = 1500
$ws.Range("M21").Value = 1400
$ws.Range("P21").Value = 1400
$ws.Range("D22").Value = 44537
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = 850
$ws.Range("P22").Value = 850
$ws.Range("D23").Value = 44504
$ws.Range("J23").Value = 500
$ws.Range("K23").Value = 900
$ws.Range("L23").Value = 1000
$ws.Range("M23").Value = 950
$ws.Range("P23").Value = 950
$ws.Range("D24").Value = 44518
$ws.Range("J24").Value = 400
$ws.Range("D25").Value = 44517
$ws.Range("K25").Value = 800
$ws.Range("L25").Value = 900
$ws.Range("M25").Value = 850
$ws.Range("P25").Value = 850
$ws.Range("D26").Value = 44503
$ws.Range("J26").Value = 400
$ws.Range("K26").Value = 900
$ws.Range("L26").Value = 1000
$ws.Range("M26").Value = 950
$ws.Range("P26").Value = 950
$ws.Range("D27").Value = 44832
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 2200
$ws.Range("M27").Value = 2100
$ws.Range("P27").Value = 2100
$ws.Range("D28").Value = 44845
$ws.Range("J28").Value = 400
$ws.Range("K28").Value = 1300
$ws.Range("L28").Value = 1500
$ws.Range("M28").Value = 1400
$ws.Range("P28").Value = 1400
